# Generate Report for Handback
# Swap the handback file identifiers / timestamps for the two rows that are
# tracked on the Overview / zh-cn / de-de sheets.
#
# Old row-2 id: 724b892e-5a63-44a9-8a22-a26a6d50ac82   -> New: 0fbc0e17-b28f-4404-b25c-e821d77e00f9
# Old row-3 id: 97a71edd-6aee-4803-89fa-0a1f2e37d95e    -> New: ffffea83ba13-756a-458d-9520-df0cdd760d30
# The handoff/handback hashes and timestamps are refreshed to match the new run.

$wb = $excel.ActiveWorkbook

$oldId1 = "724b892e-5a63-44a9-8a22-a26a6d50ac82"
$newId1 = "0fbc0e17-b28f-4404-b25c-e821d77e00f9"
$oldId2 = "97a71edd-6aee-4803-89fa-0a1f2e37d95e"
$newId2 = "ffffea83ba13-756a-458d-9520-df0cdd760d30"

$oldHash1 = "1f4387566ab9ddeb591e9375fe66ff1c5f05ac0e"
$newHash1 = "08daaee363b0626c1bf177751cd78aab34bc0beb"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "$newId1.md"
$ws.Range("B2").Hyperlinks.Item(1).TextToDisplay = "e2e\$newId1.md"
$ws.Range("G2").Value = "2016-08-16 16:57:34"

$ws.Range("A3").Value = "$newId2.md"
$ws.Range("B3").Hyperlinks.Item(1).TextToDisplay = "e2e\$newId2.md"
$ws.Range("G3").Value = "2016-08-16 16:57:34"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Hyperlinks.Item(1).TextToDisplay = "$newId1.md"
$ws.Range("G2").Value = "$newId1.$newHash1.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-16 16:57:29"
$ws.Range("I2").Hyperlinks.Item(1).TextToDisplay = "$newId1.md"
$ws.Range("J2").Value = "$newId1.$newHash1.zh-cn.xlf"
$ws.Range("K2").Value = "2016-08-16 16:57:46"

$ws.Range("A3").Hyperlinks.Item(1).TextToDisplay = "$newId2.md"
$ws.Range("G3").Value = "$newId1.$newHash1.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-16 16:57:29"
$ws.Range("I3").Hyperlinks.Item(1).TextToDisplay = "$newId2.md"
$ws.Range("J3").Value = "$newId1.$newHash1.zh-cn.xlf"
$ws.Range("K3").Value = "2016-08-16 16:57:46"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Hyperlinks.Item(1).TextToDisplay = "$newId1.md"
$ws.Range("G2").Value = "$newId1.$newHash1.de-de.xlf"
$ws.Range("H2").Value = "2016-08-16 16:57:34"
$ws.Range("I2").Hyperlinks.Item(1).TextToDisplay = "$newId1.md"
$ws.Range("J2").Value = "$newId1.$newHash1.de-de.xlf"
$ws.Range("K2").Value = "2016-08-16 16:57:53"

$ws.Range("A3").Hyperlinks.Item(1).TextToDisplay = "$newId2.md"
$ws.Range("G3").Value = "$newId1.$newHash1.de-de.xlf"
$ws.Range("H3").Value = "2016-08-16 16:57:34"
$ws.Range("I3").Hyperlinks.Item(1).TextToDisplay = "$newId2.md"
$ws.Range("J3").Value = "$newId1.$newHash1.de-de.xlf"
$ws.Range("K3").Value = "2016-08-16 16:57:53"
